$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append a new row (row 5) with the latest tracker data entry.
$ws.Range("A5").Value = "G3"
$ws.Range("B5").Value = "Test2"
$ws.Range("C5").Value = 45860
$ws.Range("C5").NumberFormat = "YYYY-MM-DD"
$ws.Range("D5").Value = 1
$ws.Range("E5").Value = 0
$ws.Range("F5").Value = 0
